# EIA Table 4.3 monthly update: November 2016 data added.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update title / period text (shared strings used by A1 title and the
#        "Rolling 12 Months Ending in ..." caption) ---------------------------
$ws.Range("A1").Value = "Table 4.3. Receipts, Average Cost, and Quality of Fossil Fuels: Independent Power Producers, 2006 - November 2016"

# --- 2. Insert the new "November" data row above the old "Year to Date" ------
#        block (old row 53), pushing everything below it down by one row, and
#        copy the number formatting down from the row above (row 52) so the
#        new row matches the existing month rows exactly. ---------------------
$ws.Rows.Item(53).Insert()
$ws.Range("A52:M52").Copy()
$ws.Range("A53:M53").PasteSpecial(-4122)

$ws.Range("A53").Value = "November"
$ws.Range("B53").Value = 270974
$ws.Range("C53").Value = 14586
$ws.Range("D53").Value = 1.92
$ws.Range("E53").Value = 35.71
$ws.Range("F53").Value = 1.54
$ws.Range("G53").Value = 114.5
$ws.Range("H53").Value = 2625
$ws.Range("I53").Value = 441
$ws.Range("J53").Value = 10.76
$ws.Range("K53").Value = 64.02
$ws.Range("L53").Value = 0.47
$ws.Range("M53").Value = 115.1

# --- 3. Refresh the "Year to Date" annual totals (now rows 55-57) -----------
$ws.Range("B55").Value = 3875808
$ws.Range("C55").Value = 206753
$ws.Range("D55").Value = 2.26
$ws.Range("E55").Value = 42.33
$ws.Range("F55").Value = 1.61
$ws.Range("G55").Value = 99
$ws.Range("H55").Value = 66408
$ws.Range("I55").Value = 11076
$ws.Range("J55").Value = 20.27
$ws.Range("K55").Value = 121.73
$ws.Range("L55").Value = 0.45
$ws.Range("M55").Value = 97.1

$ws.Range("B56").Value = 3453390
$ws.Range("C56").Value = 184146
$ws.Range("D56").Value = 2.11
$ws.Range("E56").Value = 39.6
$ws.Range("F56").Value = 1.66
$ws.Range("G56").Value = 99.2
$ws.Range("H56").Value = 52848
$ws.Range("I56").Value = 8780
$ws.Range("J56").Value = 11.78
$ws.Range("K56").Value = 70.98
$ws.Range("L56").Value = 0.46
$ws.Range("M56").Value = 86.2

$ws.Range("B57").Value = 2676066
$ws.Range("C57").Value = 142397
$ws.Range("D57").Value = 1.94
$ws.Range("E57").Value = 36.41
$ws.Range("F57").Value = 1.76
$ws.Range("G57").Value = 88.7
$ws.Range("H57").Value = 23506
$ws.Range("I57").Value = 3994
$ws.Range("J57").Value = 9.85
$ws.Range("K57").Value = 58.01
$ws.Range("L57").Value = 0.47
$ws.Range("M57").Value = 75.6

# --- 4. Update the "Rolling 12 Months Ending in ..." caption ----------------
$ws.Range("A58").Value = "Rolling 12 Months Ending in November"

# --- 5. Refresh the rolling-12-month totals (now rows 59-60) ----------------
$ws.Range("B59").Value = 3821531
$ws.Range("C59").Value = 203993
$ws.Range("D59").Value = 2.12
$ws.Range("E59").Value = 39.72
$ws.Range("F59").Value = 1.65
$ws.Range("G59").Value = 100.3
$ws.Range("H59").Value = 58213
$ws.Range("I59").Value = 9684
$ws.Range("J59").Value = 12.16
$ws.Range("K59").Value = 73.13
$ws.Range("L59").Value = 0.47
$ws.Range("M59").Value = 91.1

$ws.Range("B60").Value = 2954185
$ws.Range("C60").Value = 157232
$ws.Range("D60").Value = 1.94
$ws.Range("E60").Value = 36.45
$ws.Range("F60").Value = 1.75
$ws.Range("G60").Value = 91.1
$ws.Range("H60").Value = 25907
$ws.Range("I60").Value = 4403
$ws.Range("L60").Value = 0.47
$ws.Range("M60").Value = 76.9
# J60/K60 stay "W" (withheld) - already shifted down from the old row 59.
